$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 44: new audio entry (Jacco18 - access denied buzz / "No lab card") ---
$ws.Range("D44").Value = "https://freesound.org/people/Jacco18/sounds/419023/"
$ws.Range("B44").Value = "No lab card"
$ws.Range("A44").Value = "419023__jacco18__acess-denied-buzz.mp3"
$ws.Range("C44").Value = "https://freesound.org/"
$ws.Range("E44").Value = "https://creativecommons.org/publicdomain/zero/1.0/"

$ws.Hyperlinks.Add($ws.Range("D44"), "https://freesound.org/people/Jacco18/sounds/419023/")
$ws.Hyperlinks.Add($ws.Range("C44"), "https://freesound.org/")
$ws.Hyperlinks.Add($ws.Range("E44"), "https://creativecommons.org/publicdomain/zero/1.0/")

# --- Row 45: new audio entry (stk13 - jungle ninja / "Ambient music") ---
$ws.Range("D45").Value = "https://freesound.org/people/stk13/sounds/121980/"
$ws.Range("B45").Value = "Ambient music"
$ws.Range("A45").Value = "121980__stk13__jungle-ninja.wav"
$ws.Range("C45").Value = "https://freesound.org/"
$ws.Range("E45").Value = "https://creativecommons.org/publicdomain/zero/1.0/"

$ws.Hyperlinks.Add($ws.Range("D45"), "https://freesound.org/people/stk13/sounds/121980/")
$ws.Hyperlinks.Add($ws.Range("C45"), "https://freesound.org/")
$ws.Hyperlinks.Add($ws.Range("E45"), "https://creativecommons.org/publicdomain/zero/1.0/")

# Re-apply the same "Website / Direct link / License" formatting used by the
# rows above (Hyperlinks.Add mangles the cell style), so C44:E45 end up with
# the exact same style as the other hyperlink cells in the table.
$ws.Range("C43:E43").Copy()
$ws.Range("C44:E45").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Move the active selection, matching the saved workbook state.
$ws.Range("C31").Select()
